$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44162
$ws.Cells.Item(2, 11).Value = 'Castle Brite'
$ws.Cells.Item(2, 12).Value = 'Tercera'
$ws.Cells.Item(2, 13).Value = 500
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 16000
$ws.Cells.Item(2, 16).Value = 15500
$ws.Cells.Item(2, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 1033
$ws.Cells.Item(2, 20).Value = 15

$ws.Cells.Item(3, 4).Value = 44166
$ws.Cells.Item(3, 11).Value = 'Castle Brite'
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 600
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 17000
$ws.Cells.Item(3, 16).Value = 16500
$ws.Cells.Item(3, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 1100
$ws.Cells.Item(3, 20).Value = 15

$ws.Cells.Item(4, 4).Value = 44176
$ws.Cells.Item(4, 11).Value = 'Castle Brite'
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 13).Value = 500
$ws.Cells.Item(4, 14).Value = 15000
$ws.Cells.Item(4, 15).Value = 16000
$ws.Cells.Item(4, 16).Value = 15500
$ws.Cells.Item(4, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 1033
$ws.Cells.Item(4, 20).Value = 15

$ws.Cells.Item(5, 4).Value = 44544
$ws.Cells.Item(5, 11).Value = 'Castle Brite'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 600
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 19000
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(5, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value = 1056
$ws.Cells.Item(5, 20).Value = 18

$ws.Cells.Item(6, 4).Value = 44544
$ws.Cells.Item(6, 11).Value = 'Castle Brite'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 16000
$ws.Cells.Item(6, 15).Value = 16000
$ws.Cells.Item(6, 16).Value = 16000
$ws.Cells.Item(6, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(6, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 19).Value = 889
$ws.Cells.Item(6, 20).Value = 18

$ws.Cells.Item(7, 4).Value = 44565
$ws.Cells.Item(7, 11).Value = 'Castle Brite'
$ws.Cells.Item(7, 12).Value = 'Especial'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 20000
$ws.Cells.Item(7, 15).Value = 20000
$ws.Cells.Item(7, 16).Value = 20000
$ws.Cells.Item(7, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 1111
$ws.Cells.Item(7, 20).Value = 18

$ws.Cells.Item(8, 4).Value = 44565
$ws.Cells.Item(8, 11).Value = 'Castle Brite'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(8, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 18

$ws.Cells.Item(9, 4).Value = 44565
$ws.Cells.Item(9, 11).Value = 'Castle Brite'
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 16000
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(9, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 19).Value = 889
$ws.Cells.Item(9, 20).Value = 18

$ws.Cells.Item(10, 4).Value = 44904
$ws.Cells.Item(10, 11).Value = 'Castle Brite'
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 600
$ws.Cells.Item(10, 14).Value = 21000
$ws.Cells.Item(10, 15).Value = 22000
$ws.Cells.Item(10, 16).Value = 21500
$ws.Cells.Item(10, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 1344
$ws.Cells.Item(10, 20).Value = 16

$ws.Cells.Item(11, 4).Value = 44904
$ws.Cells.Item(11, 11).Value = 'Castle Brite'
$ws.Cells.Item(11, 12).Value = 'Segunda'
$ws.Cells.Item(11, 13).Value = 300
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 18000
$ws.Cells.Item(11, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(11, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 19).Value = 1125
$ws.Cells.Item(11, 20).Value = 16

$ws.Cells.Item(12, 4).Value = 44901
$ws.Cells.Item(12, 11).Value = 'Castle Brite'
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 600
$ws.Cells.Item(12, 14).Value = 21000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 21500
$ws.Cells.Item(12, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(12, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 19).Value = 1344
$ws.Cells.Item(12, 20).Value = 16

$ws.Cells.Item(13, 4).Value = 44901
$ws.Cells.Item(13, 11).Value = 'Castle Brite'
$ws.Cells.Item(13, 12).Value = 'Segunda'
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 18000
$ws.Cells.Item(13, 15).Value = 18000
$ws.Cells.Item(13, 16).Value = 18000
$ws.Cells.Item(13, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(13, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(13, 19).Value = 1125
$ws.Cells.Item(13, 20).Value = 16

$ws.Cells.Item(14, 4).Value = 44917
$ws.Cells.Item(14, 11).Value = 'Castle Brite'
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 300
$ws.Cells.Item(14, 14).Value = 19000
$ws.Cells.Item(14, 15).Value = 20000
$ws.Cells.Item(14, 16).Value = 19500
$ws.Cells.Item(14, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(14, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(14, 19).Value = 1219
$ws.Cells.Item(14, 20).Value = 16

$ws.Cells.Item(15, 4).Value = 44897
$ws.Cells.Item(15, 11).Value = 'Castle Brite'
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 600
$ws.Cells.Item(15, 14).Value = 23000
$ws.Cells.Item(15, 15).Value = 24000
$ws.Cells.Item(15, 16).Value = 23500
$ws.Cells.Item(15, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(15, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 19).Value = 1469
$ws.Cells.Item(15, 20).Value = 16

$ws.Cells.Item(16, 4).Value = 44897
$ws.Cells.Item(16, 11).Value = 'Castle Brite'
$ws.Cells.Item(16, 12).Value = 'Segunda'
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 19000
$ws.Cells.Item(16, 15).Value = 19000
$ws.Cells.Item(16, 16).Value = 19000
$ws.Cells.Item(16, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(16, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 19).Value = 1188
$ws.Cells.Item(16, 20).Value = 16

$ws.Cells.Item(17, 4).Value = 44943
$ws.Cells.Item(17, 11).Value = 'Modesto'
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 400
$ws.Cells.Item(17, 14).Value = 20000
$ws.Cells.Item(17, 15).Value = 21000
$ws.Cells.Item(17, 16).Value = 20500
$ws.Cells.Item(17, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(17, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(17, 19).Value = 1281
$ws.Cells.Item(17, 20).Value = 16

$ws.Cells.Item(18, 4).Value = 44547
$ws.Cells.Item(18, 11).Value = 'Castle Brite'
$ws.Cells.Item(18, 12).Value = 'Especial'
$ws.Cells.Item(18, 13).Value = 350
$ws.Cells.Item(18, 14).Value = 20000
$ws.Cells.Item(18, 15).Value = 20000
$ws.Cells.Item(18, 16).Value = 20000
$ws.Cells.Item(18, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(18, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(18, 19).Value = 1111
$ws.Cells.Item(18, 20).Value = 18

$ws.Cells.Item(19, 4).Value = 44547
$ws.Cells.Item(19, 11).Value = 'Castle Brite'
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 350
$ws.Cells.Item(19, 14).Value = 18000
$ws.Cells.Item(19, 15).Value = 18000
$ws.Cells.Item(19, 16).Value = 18000
$ws.Cells.Item(19, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(19, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 19).Value = 1000
$ws.Cells.Item(19, 20).Value = 18

$ws.Cells.Item(20, 4).Value = 44547
$ws.Cells.Item(20, 11).Value = 'Castle Brite'
$ws.Cells.Item(20, 12).Value = 'Segunda'
$ws.Cells.Item(20, 13).Value = 350
$ws.Cells.Item(20, 14).Value = 16000
$ws.Cells.Item(20, 15).Value = 16000
$ws.Cells.Item(20, 16).Value = 16000
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(20, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(20, 19).Value = 889
$ws.Cells.Item(20, 20).Value = 18

$ws.Cells.Item(21, 4).Value = 44915
$ws.Cells.Item(21, 11).Value = 'Castle Brite'
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 800
$ws.Cells.Item(21, 14).Value = 19000
$ws.Cells.Item(21, 15).Value = 20000
$ws.Cells.Item(21, 16).Value = 19500
$ws.Cells.Item(21, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(21, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 19).Value = 1219
$ws.Cells.Item(21, 20).Value = 16

$ws.Cells.Item(22, 4).Value = 44169
$ws.Cells.Item(22, 11).Value = 'Castle Brite'
$ws.Cells.Item(22, 12).Value = 'Segunda'
$ws.Cells.Item(22, 13).Value = 500
$ws.Cells.Item(22, 14).Value = 15000
$ws.Cells.Item(22, 15).Value = 16000
$ws.Cells.Item(22, 16).Value = 15500
$ws.Cells.Item(22, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(22, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(22, 19).Value = 1033
$ws.Cells.Item(22, 20).Value = 15

$ws.Cells.Item(23, 4).Value = 44925
$ws.Cells.Item(23, 11).Value = 'Castle Brite'
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 600
$ws.Cells.Item(23, 14).Value = 19000
$ws.Cells.Item(23, 15).Value = 20000
$ws.Cells.Item(23, 16).Value = 19500
$ws.Cells.Item(23, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(23, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23, 19).Value = 1219
$ws.Cells.Item(23, 20).Value = 16

$ws.Cells.Item(24, 4).Value = 44925
$ws.Cells.Item(24, 11).Value = 'Castle Brite'
$ws.Cells.Item(24, 12).Value = 'Segunda'
$ws.Cells.Item(24, 13).Value = 300
$ws.Cells.Item(24, 14).Value = 15000
$ws.Cells.Item(24, 15).Value = 15000
$ws.Cells.Item(24, 16).Value = 15000
$ws.Cells.Item(24, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(24, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(24, 19).Value = 938
$ws.Cells.Item(24, 20).Value = 16

$ws.Cells.Item(25, 4).Value = 44533
$ws.Cells.Item(25, 11).Value = 'Castle Brite'
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 350
$ws.Cells.Item(25, 14).Value = 24000
$ws.Cells.Item(25, 15).Value = 24000
$ws.Cells.Item(25, 16).Value = 24000
$ws.Cells.Item(25, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(25, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(25, 19).Value = 1333
$ws.Cells.Item(25, 20).Value = 18

$ws.Cells.Item(26, 4).Value = 44533
$ws.Cells.Item(26, 11).Value = 'Castle Brite'
$ws.Cells.Item(26, 12).Value = 'Segunda'
$ws.Cells.Item(26, 13).Value = 350
$ws.Cells.Item(26, 14).Value = 20000
$ws.Cells.Item(26, 15).Value = 20000
$ws.Cells.Item(26, 16).Value = 20000
$ws.Cells.Item(26, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(26, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(26, 19).Value = 1111
$ws.Cells.Item(26, 20).Value = 18

$ws.Cells.Item(27, 4).Value = 44533
$ws.Cells.Item(27, 11).Value = 'Castle Brite'
$ws.Cells.Item(27, 12).Value = 'Tercera'
$ws.Cells.Item(27, 13).Value = 350
$ws.Cells.Item(27, 14).Value = 17000
$ws.Cells.Item(27, 15).Value = 17000
$ws.Cells.Item(27, 16).Value = 17000
$ws.Cells.Item(27, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(27, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(27, 19).Value = 944
$ws.Cells.Item(27, 20).Value = 18

$ws.Cells.Item(28, 4).Value = 44579
$ws.Cells.Item(28, 11).Value = 'Modesto'
$ws.Cells.Item(28, 12).Value = 'Especial'
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 21000
$ws.Cells.Item(28, 15).Value = 21000
$ws.Cells.Item(28, 16).Value = 21000
$ws.Cells.Item(28, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(28, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(28, 19).Value = 1167
$ws.Cells.Item(28, 20).Value = 18

$ws.Cells.Item(29, 4).Value = 44579
$ws.Cells.Item(29, 11).Value = 'Modesto'
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 18000
$ws.Cells.Item(29, 15).Value = 18000
$ws.Cells.Item(29, 16).Value = 18000
$ws.Cells.Item(29, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(29, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 19).Value = 1000
$ws.Cells.Item(29, 20).Value = 18

$ws.Cells.Item(30, 4).Value = 44579
$ws.Cells.Item(30, 11).Value = 'Modesto'
$ws.Cells.Item(30, 12).Value = 'Segunda'
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(30, 14).Value = 16000
$ws.Cells.Item(30, 15).Value = 16000
$ws.Cells.Item(30, 16).Value = 16000
$ws.Cells.Item(30, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(30, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(30, 19).Value = 889
$ws.Cells.Item(30, 20).Value = 18

$ws.Cells.Item(31, 4).Value = 44572
$ws.Cells.Item(31, 11).Value = 'Modesto'
$ws.Cells.Item(31, 12).Value = 'Especial'
$ws.Cells.Item(31, 13).Value = 150
$ws.Cells.Item(31, 14).Value = 21000
$ws.Cells.Item(31, 15).Value = 21000
$ws.Cells.Item(31, 16).Value = 21000
$ws.Cells.Item(31, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(31, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(31, 19).Value = 1167
$ws.Cells.Item(31, 20).Value = 18

$ws.Cells.Item(32, 4).Value = 44572
$ws.Cells.Item(32, 11).Value = 'Modesto'
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 150
$ws.Cells.Item(32, 14).Value = 18000
$ws.Cells.Item(32, 15).Value = 18000
$ws.Cells.Item(32, 16).Value = 18000
$ws.Cells.Item(32, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(32, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(32, 19).Value = 1000
$ws.Cells.Item(32, 20).Value = 18

$ws.Cells.Item(33, 4).Value = 44572
$ws.Cells.Item(33, 11).Value = 'Modesto'
$ws.Cells.Item(33, 12).Value = 'Segunda'
$ws.Cells.Item(33, 13).Value = 150
$ws.Cells.Item(33, 14).Value = 16000
$ws.Cells.Item(33, 15).Value = 16000
$ws.Cells.Item(33, 16).Value = 16000
$ws.Cells.Item(33, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(33, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(33, 19).Value = 889
$ws.Cells.Item(33, 20).Value = 18

$ws.Cells.Item(34, 4).Value = 44911
$ws.Cells.Item(34, 11).Value = 'Castle Brite'
$ws.Cells.Item(34, 12).Value = 'Primera'
$ws.Cells.Item(34, 13).Value = 600
$ws.Cells.Item(34, 14).Value = 19000
$ws.Cells.Item(34, 15).Value = 20000
$ws.Cells.Item(34, 16).Value = 19500
$ws.Cells.Item(34, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(34, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(34, 19).Value = 1219
$ws.Cells.Item(34, 20).Value = 16

$ws.Cells.Item(35, 4).Value = 44924
$ws.Cells.Item(35, 11).Value = 'Castle Brite'
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 400
$ws.Cells.Item(35, 14).Value = 19000
$ws.Cells.Item(35, 15).Value = 20000
$ws.Cells.Item(35, 16).Value = 19500
$ws.Cells.Item(35, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(35, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(35, 19).Value = 1219
$ws.Cells.Item(35, 20).Value = 16

$ws.Cells.Item(36, 4).Value = 44924
$ws.Cells.Item(36, 11).Value = 'Castle Brite'
$ws.Cells.Item(36, 12).Value = 'Segunda'
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 15000
$ws.Cells.Item(36, 15).Value = 15000
$ws.Cells.Item(36, 16).Value = 15000
$ws.Cells.Item(36, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(36, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(36, 19).Value = 938
$ws.Cells.Item(36, 20).Value = 16

$ws.Cells.Item(37, 4).Value = 44523
$ws.Cells.Item(37, 11).Value = 'Castle Brite'
$ws.Cells.Item(37, 12).Value = 'Segunda'
$ws.Cells.Item(37, 13).Value = 500
$ws.Cells.Item(37, 14).Value = 28000
$ws.Cells.Item(37, 15).Value = 28500
$ws.Cells.Item(37, 16).Value = 28250
$ws.Cells.Item(37, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(37, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(37, 19).Value = 1569
$ws.Cells.Item(37, 20).Value = 18

$ws.Cells.Item(38, 4).Value = 44918
$ws.Cells.Item(38, 11).Value = 'Castle Brite'
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 600
$ws.Cells.Item(38, 14).Value = 19000
$ws.Cells.Item(38, 15).Value = 20000
$ws.Cells.Item(38, 16).Value = 19500
$ws.Cells.Item(38, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(38, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(38, 19).Value = 1219
$ws.Cells.Item(38, 20).Value = 16

$ws.Cells.Item(39, 4).Value = 44551
$ws.Cells.Item(39, 11).Value = 'Castle Brite'
$ws.Cells.Item(39, 12).Value = 'Especial'
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 20000
$ws.Cells.Item(39, 15).Value = 20000
$ws.Cells.Item(39, 16).Value = 20000
$ws.Cells.Item(39, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(39, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(39, 19).Value = 1111
$ws.Cells.Item(39, 20).Value = 18

$ws.Cells.Item(40, 4).Value = 44551
$ws.Cells.Item(40, 11).Value = 'Castle Brite'
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 200
$ws.Cells.Item(40, 14).Value = 18000
$ws.Cells.Item(40, 15).Value = 18000
$ws.Cells.Item(40, 16).Value = 18000
$ws.Cells.Item(40, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(40, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(40, 19).Value = 1000
$ws.Cells.Item(40, 20).Value = 18

$ws.Cells.Item(41, 4).Value = 44551
$ws.Cells.Item(41, 11).Value = 'Castle Brite'
$ws.Cells.Item(41, 12).Value = 'Segunda'
$ws.Cells.Item(41, 13).Value = 200
$ws.Cells.Item(41, 14).Value = 16000
$ws.Cells.Item(41, 15).Value = 16000
$ws.Cells.Item(41, 16).Value = 16000
$ws.Cells.Item(41, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(41, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(41, 19).Value = 889
$ws.Cells.Item(41, 20).Value = 18

$ws.Cells.Item(42, 4).Value = 44914
$ws.Cells.Item(42, 11).Value = 'Castle Brite'
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 400
$ws.Cells.Item(42, 14).Value = 19000
$ws.Cells.Item(42, 15).Value = 20000
$ws.Cells.Item(42, 16).Value = 19500
$ws.Cells.Item(42, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(42, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(42, 19).Value = 1219
$ws.Cells.Item(42, 20).Value = 16

$ws.Cells.Item(43, 4).Value = 44159
$ws.Cells.Item(43, 11).Value = 'Castle Brite'
$ws.Cells.Item(43, 12).Value = 'Tercera'
$ws.Cells.Item(43, 13).Value = 400
$ws.Cells.Item(43, 14).Value = 15500
$ws.Cells.Item(43, 15).Value = 16000
$ws.Cells.Item(43, 16).Value = 15750
$ws.Cells.Item(43, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(43, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(43, 19).Value = 1050
$ws.Cells.Item(43, 20).Value = 15

$ws.Cells.Item(44, 4).Value = 44900
$ws.Cells.Item(44, 11).Value = 'Castle Brite'
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 200
$ws.Cells.Item(44, 14).Value = 23000
$ws.Cells.Item(44, 15).Value = 24000
$ws.Cells.Item(44, 16).Value = 23500
$ws.Cells.Item(44, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(44, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(44, 19).Value = 1469
$ws.Cells.Item(44, 20).Value = 16

$ws.Cells.Item(45, 4).Value = 44900
$ws.Cells.Item(45, 11).Value = 'Castle Brite'
$ws.Cells.Item(45, 12).Value = 'Segunda'
$ws.Cells.Item(45, 13).Value = 100
$ws.Cells.Item(45, 14).Value = 19000
$ws.Cells.Item(45, 15).Value = 19000
$ws.Cells.Item(45, 16).Value = 19000
$ws.Cells.Item(45, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(45, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(45, 19).Value = 1188
$ws.Cells.Item(45, 20).Value = 16

$ws.Cells.Item(46, 4).Value = 44946
$ws.Cells.Item(46, 11).Value = 'Modesto'
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 600
$ws.Cells.Item(46, 14).Value = 20000
$ws.Cells.Item(46, 15).Value = 21000
$ws.Cells.Item(46, 16).Value = 20500
$ws.Cells.Item(46, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(46, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(46, 19).Value = 1281
$ws.Cells.Item(46, 20).Value = 16

$ws.Cells.Item(47, 4).Value = 44530
$ws.Cells.Item(47, 11).Value = 'Castle Brite'
$ws.Cells.Item(47, 12).Value = 'Segunda'
$ws.Cells.Item(47, 13).Value = 500
$ws.Cells.Item(47, 14).Value = 20000
$ws.Cells.Item(47, 15).Value = 21000
$ws.Cells.Item(47, 16).Value = 20500
$ws.Cells.Item(47, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(47, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(47, 19).Value = 1139
$ws.Cells.Item(47, 20).Value = 18

$ws.Cells.Item(48, 4).Value = 44895
$ws.Cells.Item(48, 11).Value = 'Castle Brite'
$ws.Cells.Item(48, 12).Value = 'Primera'
$ws.Cells.Item(48, 13).Value = 400
$ws.Cells.Item(48, 14).Value = 23000
$ws.Cells.Item(48, 15).Value = 24000
$ws.Cells.Item(48, 16).Value = 23500
$ws.Cells.Item(48, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(48, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(48, 19).Value = 1469
$ws.Cells.Item(48, 20).Value = 16

$ws.Cells.Item(49, 4).Value = 44936
$ws.Cells.Item(49, 11).Value = 'Modesto'
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 400
$ws.Cells.Item(49, 14).Value = 20000
$ws.Cells.Item(49, 15).Value = 21000
$ws.Cells.Item(49, 16).Value = 20500
$ws.Cells.Item(49, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(49, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(49, 19).Value = 1281
$ws.Cells.Item(49, 20).Value = 16

$ws.Cells.Item(50, 4).Value = 44540
$ws.Cells.Item(50, 11).Value = 'Castle Brite'
$ws.Cells.Item(50, 12).Value = 'Segunda'
$ws.Cells.Item(50, 13).Value = 600
$ws.Cells.Item(50, 14).Value = 16000
$ws.Cells.Item(50, 15).Value = 16000
$ws.Cells.Item(50, 16).Value = 16000
$ws.Cells.Item(50, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(50, 18).Value = 'Región del Maule'
$ws.Cells.Item(50, 19).Value = 889
$ws.Cells.Item(50, 20).Value = 18

$ws.Cells.Item(51, 4).Value = 44187
$ws.Cells.Item(51, 11).Value = 'Castle Brite'
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 350
$ws.Cells.Item(51, 14).Value = 16000
$ws.Cells.Item(51, 15).Value = 16000
$ws.Cells.Item(51, 16).Value = 16000
$ws.Cells.Item(51, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(51, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(51, 19).Value = 1067
$ws.Cells.Item(51, 20).Value = 15

$ws.Cells.Item(52, 4).Value = 44187
$ws.Cells.Item(52, 11).Value = 'Castle Brite'
$ws.Cells.Item(52, 12).Value = 'Segunda'
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 13000
$ws.Cells.Item(52, 15).Value = 13000
$ws.Cells.Item(52, 16).Value = 13000
$ws.Cells.Item(52, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(52, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(52, 19).Value = 867
$ws.Cells.Item(52, 20).Value = 15

$ws.Cells.Item(53, 4).Value = 44537
$ws.Cells.Item(53, 11).Value = 'Castle Brite'
$ws.Cells.Item(53, 12).Value = 'Primera'
$ws.Cells.Item(53, 13).Value = 500
$ws.Cells.Item(53, 14).Value = 20000
$ws.Cells.Item(53, 15).Value = 22000
$ws.Cells.Item(53, 16).Value = 21000
$ws.Cells.Item(53, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(53, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(53, 19).Value = 1167
$ws.Cells.Item(53, 20).Value = 18

$ws.Cells.Item(54, 4).Value = 44537
$ws.Cells.Item(54, 11).Value = 'Castle Brite'
$ws.Cells.Item(54, 12).Value = 'Segunda'
$ws.Cells.Item(54, 13).Value = 250
$ws.Cells.Item(54, 14).Value = 17000
$ws.Cells.Item(54, 15).Value = 17000
$ws.Cells.Item(54, 16).Value = 17000
$ws.Cells.Item(54, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(54, 18).Value = 'Región del Maule'
$ws.Cells.Item(54, 19).Value = 944
$ws.Cells.Item(54, 20).Value = 18

$ws.Cells.Item(55, 4).Value = 44939
$ws.Cells.Item(55, 11).Value = 'Modesto'
$ws.Cells.Item(55, 12).Value = 'Primera'
$ws.Cells.Item(55, 13).Value = 600
$ws.Cells.Item(55, 14).Value = 20000
$ws.Cells.Item(55, 15).Value = 21000
$ws.Cells.Item(55, 16).Value = 20500
$ws.Cells.Item(55, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(55, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(55, 19).Value = 1281
$ws.Cells.Item(55, 20).Value = 16

$ws.Cells.Item(56, 4).Value = 44890
$ws.Cells.Item(56, 11).Value = 'Castle Brite'
$ws.Cells.Item(56, 12).Value = 'Primera'
$ws.Cells.Item(56, 13).Value = 400
$ws.Cells.Item(56, 14).Value = 24000
$ws.Cells.Item(56, 15).Value = 25000
$ws.Cells.Item(56, 16).Value = 24500
$ws.Cells.Item(56, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(56, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(56, 19).Value = 1531
$ws.Cells.Item(56, 20).Value = 16

$ws.Cells.Item(57, 4).Value = 44907
$ws.Cells.Item(57, 11).Value = 'Castle Brite'
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 400
$ws.Cells.Item(57, 14).Value = 21000
$ws.Cells.Item(57, 15).Value = 22000
$ws.Cells.Item(57, 16).Value = 21500
$ws.Cells.Item(57, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(57, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(57, 19).Value = 1344
$ws.Cells.Item(57, 20).Value = 16

$ws.Cells.Item(58, 4).Value = 44907
$ws.Cells.Item(58, 11).Value = 'Castle Brite'
$ws.Cells.Item(58, 12).Value = 'Segunda'
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 18000
$ws.Cells.Item(58, 15).Value = 18000
$ws.Cells.Item(58, 16).Value = 18000
$ws.Cells.Item(58, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(58, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(58, 19).Value = 1125
$ws.Cells.Item(58, 20).Value = 16

$ws.Cells.Item(59, 4).Value = 44894
$ws.Cells.Item(59, 11).Value = 'Castle Brite'
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 400
$ws.Cells.Item(59, 14).Value = 23000
$ws.Cells.Item(59, 15).Value = 24000
$ws.Cells.Item(59, 16).Value = 23500
$ws.Cells.Item(59, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(59, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(59, 19).Value = 1469
$ws.Cells.Item(59, 20).Value = 16

$ws.Cells.Item(60, 4).Value = 44910
$ws.Cells.Item(60, 11).Value = 'Castle Brite'
$ws.Cells.Item(60, 12).Value = 'Primera'
$ws.Cells.Item(60, 13).Value = 600
$ws.Cells.Item(60, 14).Value = 21000
$ws.Cells.Item(60, 15).Value = 22000
$ws.Cells.Item(60, 16).Value = 21500
$ws.Cells.Item(60, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(60, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(60, 19).Value = 1344
$ws.Cells.Item(60, 20).Value = 16

$ws.Cells.Item(61, 4).Value = 44910
$ws.Cells.Item(61, 11).Value = 'Castle Brite'
$ws.Cells.Item(61, 12).Value = 'Segunda'
$ws.Cells.Item(61, 13).Value = 300
$ws.Cells.Item(61, 14).Value = 18000
$ws.Cells.Item(61, 15).Value = 18000
$ws.Cells.Item(61, 16).Value = 18000
$ws.Cells.Item(61, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(61, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(61, 19).Value = 1125
$ws.Cells.Item(61, 20).Value = 16

$ws.Cells.Item(62, 4).Value = 44553
$ws.Cells.Item(62, 11).Value = 'Castle Brite'
$ws.Cells.Item(62, 12).Value = 'Especial'
$ws.Cells.Item(62, 13).Value = 250
$ws.Cells.Item(62, 14).Value = 20000
$ws.Cells.Item(62, 15).Value = 20000
$ws.Cells.Item(62, 16).Value = 20000
$ws.Cells.Item(62, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(62, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(62, 19).Value = 1111
$ws.Cells.Item(62, 20).Value = 18

$ws.Cells.Item(63, 4).Value = 44553
$ws.Cells.Item(63, 11).Value = 'Castle Brite'
$ws.Cells.Item(63, 12).Value = 'Primera'
$ws.Cells.Item(63, 13).Value = 250
$ws.Cells.Item(63, 14).Value = 18000
$ws.Cells.Item(63, 15).Value = 18000
$ws.Cells.Item(63, 16).Value = 18000
$ws.Cells.Item(63, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(63, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(63, 19).Value = 1000
$ws.Cells.Item(63, 20).Value = 18

$ws.Cells.Item(64, 4).Value = 44553
$ws.Cells.Item(64, 11).Value = 'Castle Brite'
$ws.Cells.Item(64, 12).Value = 'Segunda'
$ws.Cells.Item(64, 13).Value = 250
$ws.Cells.Item(64, 14).Value = 16000
$ws.Cells.Item(64, 15).Value = 16000
$ws.Cells.Item(64, 16).Value = 16000
$ws.Cells.Item(64, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(64, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(64, 19).Value = 889
$ws.Cells.Item(64, 20).Value = 18

$ws.Cells.Item(65, 4).Value = 44922
$ws.Cells.Item(65, 11).Value = 'Castle Brite'
$ws.Cells.Item(65, 12).Value = 'Primera'
$ws.Cells.Item(65, 13).Value = 600
$ws.Cells.Item(65, 14).Value = 19000
$ws.Cells.Item(65, 15).Value = 20000
$ws.Cells.Item(65, 16).Value = 19500
$ws.Cells.Item(65, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(65, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(65, 19).Value = 1219
$ws.Cells.Item(65, 20).Value = 16

$ws.Cells.Item(66, 4).Value = 44922
$ws.Cells.Item(66, 11).Value = 'Castle Brite'
$ws.Cells.Item(66, 12).Value = 'Segunda'
$ws.Cells.Item(66, 13).Value = 300
$ws.Cells.Item(66, 14).Value = 15000
$ws.Cells.Item(66, 15).Value = 15000
$ws.Cells.Item(66, 16).Value = 15000
$ws.Cells.Item(66, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(66, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(66, 19).Value = 938
$ws.Cells.Item(66, 20).Value = 16

$ws.Cells.Item(67, 4).Value = 44568
$ws.Cells.Item(67, 11).Value = 'Castle Brite'
$ws.Cells.Item(67, 12).Value = 'Especial'
$ws.Cells.Item(67, 13).Value = 200
$ws.Cells.Item(67, 14).Value = 21000
$ws.Cells.Item(67, 15).Value = 21000
$ws.Cells.Item(67, 16).Value = 21000
$ws.Cells.Item(67, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(67, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(67, 19).Value = 1167
$ws.Cells.Item(67, 20).Value = 18

$ws.Cells.Item(68, 4).Value = 44568
$ws.Cells.Item(68, 11).Value = 'Castle Brite'
$ws.Cells.Item(68, 12).Value = 'Primera'
$ws.Cells.Item(68, 13).Value = 200
$ws.Cells.Item(68, 14).Value = 18000
$ws.Cells.Item(68, 15).Value = 18000
$ws.Cells.Item(68, 16).Value = 18000
$ws.Cells.Item(68, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(68, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(68, 19).Value = 1000
$ws.Cells.Item(68, 20).Value = 18

$ws.Cells.Item(69, 4).Value = 44568
$ws.Cells.Item(69, 11).Value = 'Castle Brite'
$ws.Cells.Item(69, 12).Value = 'Segunda'
$ws.Cells.Item(69, 13).Value = 200
$ws.Cells.Item(69, 14).Value = 16000
$ws.Cells.Item(69, 15).Value = 16000
$ws.Cells.Item(69, 16).Value = 16000
$ws.Cells.Item(69, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(69, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(69, 19).Value = 889
$ws.Cells.Item(69, 20).Value = 18

$ws.Cells.Item(70, 4).Value = 44908
$ws.Cells.Item(70, 11).Value = 'Castle Brite'
$ws.Cells.Item(70, 12).Value = 'Primera'
$ws.Cells.Item(70, 13).Value = 600
$ws.Cells.Item(70, 14).Value = 21000
$ws.Cells.Item(70, 15).Value = 22000
$ws.Cells.Item(70, 16).Value = 21500
$ws.Cells.Item(70, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(70, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(70, 19).Value = 1344
$ws.Cells.Item(70, 20).Value = 16

$ws.Cells.Item(71, 4).Value = 44908
$ws.Cells.Item(71, 11).Value = 'Castle Brite'
$ws.Cells.Item(71, 12).Value = 'Segunda'
$ws.Cells.Item(71, 13).Value = 300
$ws.Cells.Item(71, 14).Value = 18000
$ws.Cells.Item(71, 15).Value = 18000
$ws.Cells.Item(71, 16).Value = 18000
$ws.Cells.Item(71, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(71, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(71, 19).Value = 1125
$ws.Cells.Item(71, 20).Value = 16

$ws.Cells.Item(72, 4).Value = 44194
$ws.Cells.Item(72, 11).Value = 'Castle Brite'
$ws.Cells.Item(72, 12).Value = 'Segunda'
$ws.Cells.Item(72, 13).Value = 300
$ws.Cells.Item(72, 14).Value = 15000
$ws.Cells.Item(72, 15).Value = 16000
$ws.Cells.Item(72, 16).Value = 15500
$ws.Cells.Item(72, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(72, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(72, 19).Value = 1033
$ws.Cells.Item(72, 20).Value = 15

$ws.Cells.Item(73, 4).Value = 44575
$ws.Cells.Item(73, 11).Value = 'Modesto'
$ws.Cells.Item(73, 12).Value = 'Especial'
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 21000
$ws.Cells.Item(73, 15).Value = 21000
$ws.Cells.Item(73, 16).Value = 21000
$ws.Cells.Item(73, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(73, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(73, 19).Value = 1167
$ws.Cells.Item(73, 20).Value = 18

$ws.Cells.Item(74, 4).Value = 44575
$ws.Cells.Item(74, 11).Value = 'Modesto'
$ws.Cells.Item(74, 12).Value = 'Primera'
$ws.Cells.Item(74, 13).Value = 200
$ws.Cells.Item(74, 14).Value = 18000
$ws.Cells.Item(74, 15).Value = 18000
$ws.Cells.Item(74, 16).Value = 18000
$ws.Cells.Item(74, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(74, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(74, 19).Value = 1000
$ws.Cells.Item(74, 20).Value = 18

$ws.Cells.Item(75, 4).Value = 44575
$ws.Cells.Item(75, 11).Value = 'Modesto'
$ws.Cells.Item(75, 12).Value = 'Segunda'
$ws.Cells.Item(75, 13).Value = 200
$ws.Cells.Item(75, 14).Value = 16000
$ws.Cells.Item(75, 15).Value = 16000
$ws.Cells.Item(75, 16).Value = 16000
$ws.Cells.Item(75, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(75, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(75, 19).Value = 889
$ws.Cells.Item(75, 20).Value = 18
